$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Julio de 2020 a las 17:20"

# Row 4
$ws.Range("B4").Value = 4583595
$ws.Range("C4").Value = 15558
$ws.Range("D4").Value = 2246902
$ws.Range("E4").Value = 2182507
$ws.Range("G4").Value = 346
$ws.Range("H4").Value = 154186

# Row 5
$ws.Range("B5").Value = 2556207
$ws.Range("C5").Value = 689
$ws.Range("E5").Value = 678576
$ws.Range("G5").Value = 24
$ws.Range("H5").Value = 90212

# Row 6
$ws.Range("B6").Value = 1632422
$ws.Range("C6").Value = 48038
$ws.Range("E6").Value = 567641
$ws.Range("G6").Value = 709
$ws.Range("H6").Value = 35712

# Row 18
$ws.Range("D18").Value = 199796
$ws.Range("E18").Value = 11848
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 35132

# Row 21
$ws.Range("B21").Value = 209300
$ws.Range("C21").Value = 489
$ws.Range("E21").Value = 8083
$ws.Range("G21").Value = 5
$ws.Range("H21").Value = 9217

# Row 33
$ws.Range("B33").Value = 80100
$ws.Range("C33").Value = 34
$ws.Range("G33").Value = 9
$ws.Range("H33").Value = 5739

# Row 45
$ws.Range("D45").Value = 46308
$ws.Range("E45").Value = 5474

# Row 89
$ws.Range("B89").Value = 7728
$ws.Range("C89").Value = 81
$ws.Range("D89").Value = 6270
$ws.Range("E89").Value = 1415

# Row 106
$ws.Range("B106").Value = 3858
$ws.Range("C106").Value = 120
$ws.Range("D106").Value = 1760
$ws.Range("E106").Value = 1991
$ws.Range("G106").Value = 4
$ws.Range("H106").Value = 107

# Row 113
$ws.Range("E113").Value = 1467
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 25

# Row 116
$ws.Range("B116").Value = 2962
$ws.Range("C116").Value = 57
$ws.Range("D116").Value = 2721
$ws.Range("E116").Value = 202

# Row 138
$ws.Range("B138").Value = 1514
$ws.Range("C138").Value = 26
$ws.Range("D138").Value = 1187
$ws.Range("E138").Value = 277

# Row 145
$ws.Range("B145").Value = 1134
$ws.Range("C145").Value = 2
$ws.Range("D145").Value = 1028
$ws.Range("E145").Value = 37

# Row 150
$ws.Range("B150").Value = 922
$ws.Range("C150").Value = 4
$ws.Range("D150").Value = 806
$ws.Range("E150").Value = 64

# Row 160
$ws.Range("B160").Value = 604
$ws.Range("C160").Value = 28
$ws.Range("D160").Value = 144
$ws.Range("E160").Value = 447

# Row 169
$ws.Range("D169").Value = 295
$ws.Range("E169").Value = 52

$wb.Save()
